$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final table data (players, positions, teams) for rows 2..18
$data = @(
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Isaiah Collier", "PG,SG", "Utah Jazz"),
    @("De'Aaron Fox", "PG", "San Antonio Spurs"),
    @("Luka Doncic", "PG,SG", "Los Angeles Lakers"),
    @("Scottie Barnes", "PG,SG,SF,PF", "Toronto Raptors"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Yves Missi", "C", "New Orleans Pelicans"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Royce O'Neale", "SF,PF", "Phoenix Suns"),
    @("Donovan Clingan", "C", "Portland Trail Blazers"),
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("P.J. Washington", "SF,PF", "Dallas Mavericks")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
